$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2144927536231884
$ws.Range("C2").Value = 0.5304347826086957
$ws.Range("J2").Value = 0.01159420289855072
$ws.Range("P2").Value = 0.1536231884057971
$ws.Range("S2").Value = 0.08985507246376812
$ws.Range("C3").Value = 0.0267379679144385
$ws.Range("J3").Value = 0.0106951871657754
$ws.Range("P3").Value = 0.7647058823529411
$ws.Range("S3").Value = 0.1978609625668449
$ws.Range("J4").Value = 0.08771929824561403
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.1754385964912281
$ws.Range("B6").Value = 0.04424778761061947
$ws.Range("D6").Value = 0.02212389380530973
$ws.Range("F6").Value = 0.07964601769911504
$ws.Range("J6").Value = 0.2876106194690266
$ws.Range("O6").Value = 0.008849557522123894
$ws.Range("Q6").Value = 0.1769911504424779
$ws.Range("R6").Value = 0.04867256637168142
$ws.Range("S6").Value = 0.331858407079646
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.0392156862745098
$ws.Range("F7").Value = 0.0457516339869281
$ws.Range("J7").Value = 0.1437908496732026
$ws.Range("O7").Value = 0.0261437908496732
$ws.Range("Q7").Value = 0.1764705882352941
$ws.Range("R7").Value = 0.0718954248366013
$ws.Range("S7").Value = 0.3856209150326798
$ws.Range("B8").Value = 0.1245059288537549
$ws.Range("D8").Value = 0.03952569169960474
$ws.Range("F8").Value = 0.07312252964426877
$ws.Range("J8").Value = 0.1284584980237154
$ws.Range("O8").Value = 0.01383399209486166
$ws.Range("Q8").Value = 0.158102766798419
$ws.Range("R8").Value = 0.08102766798418973
$ws.Range("S8").Value = 0.3814229249011858
$ws.Range("B9").Value = 0.1118421052631579
$ws.Range("D9").Value = 0.0131578947368421
$ws.Range("F9").Value = 0.07236842105263158
$ws.Range("J9").Value = 0.07894736842105263
$ws.Range("O9").Value = 0.04605263157894737
$ws.Range("Q9").Value = 0.1644736842105263
$ws.Range("R9").Value = 0.131578947368421
$ws.Range("S9").Value = 0.3815789473684211
$ws.Range("B10").Value = 0.1292962356792144
$ws.Range("D10").Value = 0.0204582651391162
$ws.Range("F10").Value = 0.06382978723404255
$ws.Range("J10").Value = 0.1153846153846154
$ws.Range("O10").Value = 0.0220949263502455
$ws.Range("Q10").Value = 0.2004909983633388
$ws.Range("R10").Value = 0.07037643207855974
$ws.Range("S10").Value = 0.3780687397708674
$ws.Range("G11").Value = 0.1485507246376812
$ws.Range("J11").Value = 0.1123188405797101
$ws.Range("K11").Value = 0.2210144927536232
$ws.Range("L11").Value = 0.5036231884057971
$ws.Range("S11").Value = 0.01449275362318841
$ws.Range("G12").Value = 0.6595744680851063
$ws.Range("J12").Value = 0.2836879432624114
$ws.Range("L12").Value = 0.01418439716312057
$ws.Range("S12").Value = 0.0425531914893617
$ws.Range("G13").Value = 0.6285714285714286
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.08571428571428572
$ws.Range("F15").Value = 0.03626943005181347
$ws.Range("H15").Value = 0.1865284974093264
$ws.Range("I15").Value = 0.04663212435233161
$ws.Range("J15").Value = 0.3264248704663212
$ws.Range("K15").Value = 0.06217616580310881
$ws.Range("O15").Value = 0.04663212435233161
$ws.Range("S15").Value = 0.2953367875647668
$ws.Range("F16").Value = 0.03524229074889868
$ws.Range("H16").Value = 0.1938325991189427
$ws.Range("I16").Value = 0.05286343612334802
$ws.Range("J16").Value = 0.4229074889867842
$ws.Range("K16").Value = 0.118942731277533
$ws.Range("M16").Value = 0.01762114537444934
$ws.Range("N16").Value = 0.004405286343612335
$ws.Range("O16").Value = 0.05286343612334802
$ws.Range("S16").Value = 0.1013215859030837
$ws.Range("F17").Value = 0.02877697841726619
$ws.Range("H17").Value = 0.2470023980815348
$ws.Range("I17").Value = 0.0671462829736211
$ws.Range("J17").Value = 0.3836930455635492
$ws.Range("K17").Value = 0.07913669064748201
$ws.Range("M17").Value = 0.009592326139088728
$ws.Range("N17").Value = 0.002398081534772182
$ws.Range("O17").Value = 0.05755395683453238
$ws.Range("S17").Value = 0.1247002398081535
$ws.Range("F18").Value = 0.005917159763313609
$ws.Range("H18").Value = 0.2248520710059172
$ws.Range("I18").Value = 0.09467455621301775
$ws.Range("J18").Value = 0.4437869822485207
$ws.Range("K18").Value = 0.07692307692307693
$ws.Range("M18").Value = 0.01775147928994083
$ws.Range("O18").Value = 0.04142011834319527
$ws.Range("S18").Value = 0.09467455621301775
$ws.Range("F19").Value = 0.01755786113328013
$ws.Range("H19").Value = 0.2274541101356744
$ws.Range("I19").Value = 0.06863527533918595
$ws.Range("J19").Value = 0.3591380686352754
$ws.Range("K19").Value = 0.1005586592178771
$ws.Range("M19").Value = 0.01995211492418196
$ws.Range("N19").Value = 0.003192338387869114
$ws.Range("O19").Value = 0.05506783719074222
$ws.Range("S19").Value = 0.1484437350359138
